# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-56) from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C56").Value = 45224
